$d = $word.ActiveDocument

# 1. "forbruger delen" -> "forbrugerapplikationen"
$d.Content.Find.Execute("forbruger delen", $true, $false, $false, $false, $false,
                         $true, 1, $false, "forbrugerapplikationen", 2) | Out-Null

# 2. "muligheder for at begrænse" -> "mulighed for at begrænse"
$d.Content.Find.Execute("muligheder for at begrænse", $true, $false, $false, $false, $false,
                         $true, 1, $false, "mulighed for at begrænse", 2) | Out-Null

# 3. "ingrediens liste" -> "ingrediensliste"
$d.Content.Find.Execute("ingrediens liste", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ingrediensliste", 2) | Out-Null

# 4. "forretningsmanager delen" -> "forretningsmanagerapplikationen"
$d.Content.Find.Execute("forretningsmanager delen", $true, $false, $false, $false, $false,
                         $true, 1, $false, "forretningsmanagerapplikationen", 2) | Out-Null

# 5. "tilbud priserne" -> "tilbudspriserne"
$d.Content.Find.Execute("tilbud priserne", $true, $false, $false, $false, $false,
                         $true, 1, $false, "tilbudspriserne", 2) | Out-Null
